$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6096.4443
$ws.Range("I86").Value = 653.63635
$ws.Range("J86").Value = 14649.429
$ws.Range("K86").Value = 653.63635
$ws.Range("L86").Value = 14649.429
$ws.Range("M86").Value = 469.36365
$ws.Range("N86").Value = -16895.429

$ws.Range("H89").Value = 6096.4443
$ws.Range("I89").Value = 653.63635
$ws.Range("J89").Value = 14649.429
$ws.Range("K89").Value = 3268.18175
$ws.Range("L89").Value = 73247.145
$ws.Range("M89").Value = 2347.81825
$ws.Range("N89").Value = -84479.145

$ws.Range("H98").Value = 1833.3334
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 1000
$ws.Range("M98").Value = 498

$ws.Range("H122").Value = 1833.3334
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H132").Value = 2348.925
$ws.Range("I132").Value = 2761.0322
$ws.Range("K132").Value = 8283.096600000001
$ws.Range("M132").Value = -5753.096600000001

$ws.Range("H137").Value = 1348.9615
$ws.Range("I137").Value = 1342.619
$ws.Range("J137").Value = 1375.6
$ws.Range("K137").Value = 4027.857
$ws.Range("L137").Value = 4126.799999999999
$ws.Range("M137").Value = -1477.857
$ws.Range("N137").Value = -9226.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3377.3416
$ws.Range("I32").Value = 3056.7605
$ws.Range("J32").Value = 5446.5454
$ws.Range("K32").Value = 3056.7605
$ws.Range("L32").Value = 5446.5454
$ws.Range("M32").Value = -2769.7605
$ws.Range("N32").Value = -6020.5454

$ws.Range("H61").Value = 296129.78
$ws.Range("I61").Value = 360860.38
$ws.Range("J61").Value = 1899.909
$ws.Range("K61").Value = 360860.38
$ws.Range("L61").Value = 1899.909
$ws.Range("M61").Value = -360648.38
$ws.Range("N61").Value = -2323.909

$ws.Range("H74").Value = 35716416
$ws.Range("I74").Value = 43480584
$ws.Range("J74").Value = 1251.8
$ws.Range("K74").Value = 43480584
$ws.Range("L74").Value = 1251.8
$ws.Range("M74").Value = -43479710
$ws.Range("N74").Value = -2999.8

$ws.Range("H77").Value = 35716416
$ws.Range("I77").Value = 43480584
$ws.Range("J77").Value = 1251.8
$ws.Range("K77").Value = 217402920
$ws.Range("L77").Value = 6259
$ws.Range("M77").Value = -217398552
$ws.Range("N77").Value = -14995

$ws.Range("H132").Value = 16640.47
$ws.Range("I132").Value = 1861.6296
$ws.Range("J132").Value = 73644.57000000001
$ws.Range("K132").Value = 5584.8888
$ws.Range("L132").Value = 220933.71
$ws.Range("M132").Value = -3054.8888
$ws.Range("N132").Value = -225993.71

$ws.Range("H136").Value = 296129.78
$ws.Range("I136").Value = 360860.38
$ws.Range("J136").Value = 1899.909
$ws.Range("K136").Value = 1082581.14
$ws.Range("L136").Value = 5699.727000000001
$ws.Range("M136").Value = -1080031.14
$ws.Range("N136").Value = -10799.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1540.8043
$ws.Range("I86").Value = 1317.1428
$ws.Range("J86").Value = 1888.7222
$ws.Range("K86").Value = 1317.1428
$ws.Range("L86").Value = 1888.7222
$ws.Range("M86").Value = -194.1428000000001
$ws.Range("N86").Value = -4134.7222

$ws.Range("H89").Value = 1540.8043
$ws.Range("I89").Value = 1317.1428
$ws.Range("J89").Value = 1888.7222
$ws.Range("K89").Value = 6585.714
$ws.Range("L89").Value = 9443.610999999999
$ws.Range("M89").Value = -969.7139999999999
$ws.Range("N89").Value = -20675.611

$ws.Range("H134").Value = 3555.2424
$ws.Range("I134").Value = 3859.4482
$ws.Range("J134").Value = 1349.75
$ws.Range("K134").Value = 11578.3446
$ws.Range("L134").Value = 4049.25
$ws.Range("M134").Value = -9043.3446
$ws.Range("N134").Value = -9119.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1500
$ws.Range("J4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("N4").Value = -1724

$ws.Range("H31").Value = 2980.5637
$ws.Range("I31").Value = 1787.7106
$ws.Range("J31").Value = 5646.9414
$ws.Range("K31").Value = 1787.7106
$ws.Range("L31").Value = 5646.9414
$ws.Range("M31").Value = -1492.7106
$ws.Range("N31").Value = -6236.9414

$ws.Range("H34").Value = 2980.5637
$ws.Range("I34").Value = 1787.7106
$ws.Range("J34").Value = 5646.9414
$ws.Range("K34").Value = 1787.7106
$ws.Range("L34").Value = 5646.9414
$ws.Range("M34").Value = -1585.7106
$ws.Range("N34").Value = -6050.9414

$ws.Range("H58").Value = 11480.625
$ws.Range("I58").Value = 887.94446
$ws.Range("K58").Value = 887.94446
$ws.Range("M58").Value = -684.94446

$ws.Range("H134").Value = 623.6818
$ws.Range("I134").Value = 569.79486
$ws.Range("J134").Value = 1044
$ws.Range("K134").Value = 1709.38458
$ws.Range("L134").Value = 3132
$ws.Range("M134").Value = 825.6154200000001
$ws.Range("N134").Value = -8202

$ws.Range("H136").Value = 11480.625
$ws.Range("I136").Value = 887.94446
$ws.Range("K136").Value = 2663.83338
$ws.Range("M136").Value = -113.83338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 6749.4
$ws.Range("I2").Value = 10082.1
$ws.Range("K2").Value = 60492.60000000001
$ws.Range("M2").Value = -60379.60000000001

$ws.Range("H131").Value = 783.58
$ws.Range("J131").Value = 796.4
$ws.Range("L131").Value = 2389.2
$ws.Range("N131").Value = -12469.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4667.577
$ws.Range("I7").Value = 2975.1428
$ws.Range("K7").Value = 2975.1428
$ws.Range("M7").Value = -2863.1428

$ws.Range("H61").Value = 3540.8
$ws.Range("I61").Value = 1314.0667
$ws.Range("K61").Value = 1314.0667
$ws.Range("M61").Value = -1112.0667

$ws.Range("H100").Value = 1769.9231
$ws.Range("J100").Value = 2831
$ws.Range("L100").Value = 2831
$ws.Range("N100").Value = -3913

$ws.Range("H113").Value = 3540.8
$ws.Range("I113").Value = 1314.0667
$ws.Range("K113").Value = 1314.0667
$ws.Range("M113").Value = 855.9332999999999

$ws.Range("H126").Value = 4667.577
$ws.Range("I126").Value = 2975.1428
$ws.Range("K126").Value = 8925.428400000001
$ws.Range("M126").Value = -6455.428400000001

$ws.Range("H136").Value = 829.9091
$ws.Range("I136").Value = 829.9091
$ws.Range("K136").Value = 2489.7273
$ws.Range("M136").Value = 60.27269999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2525801
$ws.Range("I107").Value = 490.7857
$ws.Range("K107").Value = 1472.3571
$ws.Range("M107").Value = 447.6428999999998
